$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing rows (col1/col3 and a3/4)
# shift down to rows 2 and 3.
$ws.Rows.Item(1).Insert()

# New header row
$ws.Range("A1").Value = "col1"
$ws.Range("B1").Value = "col2"
$ws.Range("C1").Value = "col3"

# Row 2: was the original header row (col1/col3), now becomes the first
# data row -- overwrite with the new data values.
$ws.Range("A2").Value = "a1"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3.1

# Row 3: was the original data row (a3/4), already correct in A3/B3 -- just
# fill in the new column C value.
$ws.Range("C3").Value = 5.1

# Row 4: brand new row appended at the bottom.
# A4 must be stored as TEXT "6" (not the number 6). A plain .Value="6"
# assignment gets auto-coerced to a number, so build it as a text formula
# and paste the computed value back as a literal (keeps the default style,
# no NumberFormat/style churn).
$ws.Range("E4").Formula = '="6"'
$ws.Range("E4").Copy()
$ws.Range("A4").PasteSpecial(-4163)
$ws.Range("E4").ClearContents()

$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 8.1
